$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (coins and lost), add new ratio column value
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 9
$ws.Range("H2").Value = 0.11

# Add new row 3 for "trump"
$ws.Range("A3").Value = "trump"
$ws.Range("B3").Value = "f"
$ws.Range("C3").Value = "trump.png"
$ws.Range("D3").Value = "a"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 1

# Update the selected cell as in the diff
$ws.Range("H4").Select()
